$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "Invoice Data"

# --- Header row (row 1) ---
$headers = @(
    "Invoice & Order Identification",
    "Dates",
    "Transaction Details",
    "Billing & Shipping Information",
    "Itemized Details",
    "Totals"
)

$cols = @("A", "B", "C", "D", "E", "F")

# Copy the formatting already present on A1 (bold font, border, centered/top
# aligned) onto the rest of the header row before writing the new text.
$ws.Range("A1").Copy()
for ($i = 1; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "1").PasteSpecial(-4122)
}

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "1").Value = $headers[$i]
}

# --- Data row (row 2) ---
$values = @(
    "Invoice Number: 10525RD001820847`n-----------------------------------------`nOrder Number: 1290196-056445-8121802`n-----------------------------------------`nPacket/Reference ID: 200755602",
    "Invoice Date: 31 Jan 2025`n-----------------------------------------`nOrder Date: 31 Jan 2025",
    "Nature of Transaction: In-State`n-----------------------------------------`nNature of Supply: Goods`n-----------------------------------------`nPlace of Supply: UTTAR PRADESH Nature of Supply",
    "Bill To: Not Found`n-----------------------------------------`nBill From/Ship From: Not Found`n-----------------------------------------`nGSTIN Number: O9AALCRSO22RIZN",
    "Item/Product Code: RDTPCASH101503720(R1.L0821-39)`n-----------------------------------------`nProduct Description: Red Tape Women Round Toe Lace-Ups Sneakers`n-----------------------------------------`nHSN/SAC Code: 64041990",
    "Totals: Rs 6890.00 Rs565800 _Rs0.00 Rs 1051.70 Rs94.65 Rs 94.65 Rs 1241.00"
)

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "2").Value = $values[$i]
}

# Re-measure the row so the multi-line content doesn't leave a stale,
# explicit "custom height" on the row.
$ws.Rows.Item(2).AutoFit()

# Page margins (Excel's PageSetup works in points, so convert from inches).
$ws.PageSetup.LeftMargin = 0.7 * 72
$ws.PageSetup.RightMargin = 0.7 * 72
$ws.PageSetup.TopMargin = 0.75 * 72
$ws.PageSetup.BottomMargin = 0.75 * 72
$ws.PageSetup.HeaderMargin = 0.3 * 72
$ws.PageSetup.FooterMargin = 0.3 * 72
